# Generate Report for Handback
#
# The d80c0318-0ed5-4aba-a368-0f0b4da8e66c.md file has come back from
# handback processing: its status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" on every sheet, and the per-language
# sheets (zh-cn, de-de) get their "Latest Target File" / "Latest Handback
# File" hyperlinks populated plus a fresh "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---- Overview sheet: row 5 is the d80c0318 file ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B5").Value = $handedBack
$overview.Range("C5").Value = $handedBack

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B5").Value = $handedBack
$zh.Hyperlinks.Add($zh.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ed2d48c1d9f84b502e84d7f2aebc7bb0b81cdd1e/e2e/d80c0318-0ed5-4aba-a368-0f0b4da8e66c.md", "", "", "d80c0318-0ed5-4aba-a368-0f0b4da8e66c.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0e30ea830fb34706d1080ada98a6eedc77e1730a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/d80c0318-0ed5-4aba-a368-0f0b4da8e66c.251a662edb8a031190bef07b2d1849e66c075431.zh-cn.xlf", "", "", "d80c0318-0ed5-4aba-a368-0f0b4da8e66c.251a662edb8a031190bef07b2d1849e66c075431.zh-cn.xlf") | Out-Null
$zh.Range("G5").Value = "2016-03-09 05:55:46"

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("B5").Value = $handedBack
$de.Hyperlinks.Add($de.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ed2d48c1d9f84b502e84d7f2aebc7bb0b81cdd1e/e2e/d80c0318-0ed5-4aba-a368-0f0b4da8e66c.md", "", "", "d80c0318-0ed5-4aba-a368-0f0b4da8e66c.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c741059930e2ea3e6d1048e4413e2f40d1a72a62/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/d80c0318-0ed5-4aba-a368-0f0b4da8e66c.251a662edb8a031190bef07b2d1849e66c075431.de-de.xlf", "", "", "d80c0318-0ed5-4aba-a368-0f0b4da8e66c.251a662edb8a031190bef07b2d1849e66c075431.de-de.xlf") | Out-Null
$de.Range("G5").Value = "2016-03-09 05:56:00"
